# Added New Mac-Address and Document Types
# Append 5 new rows (157-161) of device records to the
# master-reg_center_machine_device sheet, mirroring the existing data
# pattern (regcntr_id, machine_id, device_id, lang_code, is_active,
# cr_by, cr_dtimes, eff_dtimes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(10002, 10032, 3000176),
    @(10002, 10032, 3000177),
    @(10002, 10032, 3000178),
    @(10002, 10032, 3000179),
    @(10002, 10032, 3000180)
)

$startRow = 157
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
}

# Move the view/selection to reflect the newly added data (cosmetic,
# matches the author's recorded cursor position after the edit).
$ws.Range("D157").Select()
try {
    $excel.ActiveWindow.ScrollRow = 151
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Scroll-position persistence isn't critical to the data change.
}

# The author also switched the workbook to manual calculation mode.
$excel.Calculation = -4135
